# Weekly update: insert the new week's price report as a new row right
# after the existing header + first 16 data rows (i.e. at row 18), pushing
# all subsequent data rows down by one. The sheet's used range grows from
# A1:R97 to A1:R98 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; Excel shifts rows 18:97 down to 19:98
# and extends the sheet dimension automatically.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's data.
$ws.Range("A18").Value = 10
$ws.Range("B18").Value = "Vega Modelo de Temuco"
$ws.Range("C18").Value = "La Araucanía"
$ws.Range("D18").Value = 44749
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 100114002
$ws.Range("G18").Value = "Camote"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("M18").Value = 20000
$ws.Range("N18").Value = "$/malla 20 kilos"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 1000
$ws.Range("Q18").Value = 20
$ws.Range("R18").Value = "Hortaliza"
